$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number banner, and week-covering date range) ---
$ws.Range("A8").Value = "Volume 29   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/12/2022  Through  12/18/2022"

# --- Crime-complaints table updates (rows 14-30) ---
# For numeric targets we set NumberFormat first (matches the column's existing
# style: "#,##0" for counts, "#,##0.0;""-""#,##0.0" / "#,##0.00;""-""#,##0.00"
# for % columns) so the COM layer reuses the existing cellXf instead of minting
# a new one. For "N/A" targets (blank comparison => literal "0" or "***.*" text)
# we copy from the untouched template cells in row 23, which already hold those
# exact shared strings with the correct style.
$ws.Range("C23").Copy($ws.Range("C14"))
$ws.Range("C23").Copy($ws.Range("D15"))
$ws.Range("E23").Copy($ws.Range("E15"))
$ws.Range("C23").Copy($ws.Range("F15"))
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H15").Value = -100
$ws.Range("N15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N15").Value = -50
$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("C16").Value = 11
$ws.Range("E16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E16").Value = 120
$ws.Range("F16").NumberFormat = '#,##0'
$ws.Range("F16").Value = 53
$ws.Range("G16").NumberFormat = '#,##0'
$ws.Range("G16").Value = 31
$ws.Range("H16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H16").Value = 70.967741935483
$ws.Range("I16").NumberFormat = '#,##0'
$ws.Range("I16").Value = 630
$ws.Range("J16").NumberFormat = '#,##0'
$ws.Range("J16").Value = 415
$ws.Range("K16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K16").Value = 51.807228915662
$ws.Range("L16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L16").Value = 279.518072289157
$ws.Range("M16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M16").Value = 288.888888888889
$ws.Range("N16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N16").Value = -74.116680361544
$ws.Range("C17").NumberFormat = '#,##0'
$ws.Range("C17").Value = 16
$ws.Range("D17").NumberFormat = '#,##0'
$ws.Range("D17").Value = 6
$ws.Range("E17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E17").Value = 166.666666666667
$ws.Range("F17").NumberFormat = '#,##0'
$ws.Range("F17").Value = 44
$ws.Range("G17").NumberFormat = '#,##0'
$ws.Range("G17").Value = 17
$ws.Range("H17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H17").Value = 158.823529411765
$ws.Range("I17").NumberFormat = '#,##0'
$ws.Range("I17").Value = 472
$ws.Range("J17").NumberFormat = '#,##0'
$ws.Range("J17").Value = 426
$ws.Range("K17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K17").Value = 10.798122065727
$ws.Range("L17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L17").Value = 123.696682464455
$ws.Range("M17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M17").Value = 162.222222222222
$ws.Range("N17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N17").Value = -26.821705426356
$ws.Range("C18").NumberFormat = '#,##0'
$ws.Range("C18").Value = 9
$ws.Range("D18").NumberFormat = '#,##0'
$ws.Range("D18").Value = 15
$ws.Range("E18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E18").Value = -40
$ws.Range("F18").NumberFormat = '#,##0'
$ws.Range("F18").Value = 49
$ws.Range("G18").NumberFormat = '#,##0'
$ws.Range("G18").Value = 50
$ws.Range("H18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H18").Value = -2
$ws.Range("I18").NumberFormat = '#,##0'
$ws.Range("I18").Value = 642
$ws.Range("J18").NumberFormat = '#,##0'
$ws.Range("J18").Value = 435
$ws.Range("K18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K18").Value = 47.586206896551
$ws.Range("L18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L18").Value = 84.482758620689
$ws.Range("M18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M18").Value = 91.071428571428
$ws.Range("N18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N18").Value = -75.183610359489
$ws.Range("C19").NumberFormat = '#,##0'
$ws.Range("C19").Value = 61
$ws.Range("D19").NumberFormat = '#,##0'
$ws.Range("D19").Value = 41
$ws.Range("E19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E19").Value = 48.780487804878
$ws.Range("F19").NumberFormat = '#,##0'
$ws.Range("F19").Value = 256
$ws.Range("G19").NumberFormat = '#,##0'
$ws.Range("G19").Value = 164
$ws.Range("H19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H19").Value = 56.097560975609
$ws.Range("I19").NumberFormat = '#,##0'
$ws.Range("I19").Value = 2267
$ws.Range("J19").NumberFormat = '#,##0'
$ws.Range("J19").Value = 1387
$ws.Range("K19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K19").Value = 63.446286950252
$ws.Range("L19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L19").Value = 102.04991087344
$ws.Range("M19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M19").Value = 2.117117117117
$ws.Range("N19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N19").Value = -75.011022927689
$ws.Range("C20").NumberFormat = '#,##0'
$ws.Range("C20").Value = 1
$ws.Range("D20").NumberFormat = '#,##0'
$ws.Range("D20").Value = 2
$ws.Range("E20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E20").Value = -50
$ws.Range("F20").NumberFormat = '#,##0'
$ws.Range("F20").Value = 4
$ws.Range("G20").NumberFormat = '#,##0'
$ws.Range("G20").Value = 4
$ws.Range("I20").NumberFormat = '#,##0'
$ws.Range("I20").Value = 64
$ws.Range("J20").NumberFormat = '#,##0'
$ws.Range("J20").Value = 54
$ws.Range("K20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K20").Value = 18.518518518518
$ws.Range("L20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L20").Value = 36.170212765957
$ws.Range("M20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M20").Value = 190.909090909091
$ws.Range("N20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N20").Value = -81.395348837209
$ws.Range("C21").NumberFormat = '#,##0'
$ws.Range("C21").Value = 98
$ws.Range("D21").NumberFormat = '#,##0'
$ws.Range("D21").Value = 69
$ws.Range("E21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("E21").Value = 42.028985507246
$ws.Range("F21").NumberFormat = '#,##0'
$ws.Range("F21").Value = 407
$ws.Range("G21").NumberFormat = '#,##0'
$ws.Range("G21").Value = 270
$ws.Range("H21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("H21").Value = 50.740740740740
$ws.Range("I21").NumberFormat = '#,##0'
$ws.Range("I21").Value = 4102
$ws.Range("J21").NumberFormat = '#,##0'
$ws.Range("J21").Value = 2739
$ws.Range("K21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("K21").Value = 49.762687112084
$ws.Range("L21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("L21").Value = 113.534617386778
$ws.Range("M21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("M21").Value = 39.666326183180
$ws.Range("N21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("N21").Value = -72.891884747554
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("C22").Value = 4
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("D22").Value = 4
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E22").Value = 0
$ws.Range("F22").NumberFormat = '#,##0'
$ws.Range("F22").Value = 10
$ws.Range("G22").NumberFormat = '#,##0'
$ws.Range("G22").Value = 19
$ws.Range("H22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H22").Value = -47.368421052631
$ws.Range("I22").NumberFormat = '#,##0'
$ws.Range("I22").Value = 179
$ws.Range("J22").NumberFormat = '#,##0'
$ws.Range("J22").Value = 154
$ws.Range("K22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K22").Value = 16.233766233766
$ws.Range("L22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L22").Value = 37.692307692307
$ws.Range("M22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M22").Value = 23.448275862069
$ws.Range("C24").NumberFormat = '#,##0'
$ws.Range("C24").Value = 84
$ws.Range("D24").NumberFormat = '#,##0'
$ws.Range("D24").Value = 72
$ws.Range("E24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E24").Value = 16.666666666666
$ws.Range("F24").NumberFormat = '#,##0'
$ws.Range("F24").Value = 282
$ws.Range("G24").NumberFormat = '#,##0'
$ws.Range("G24").Value = 247
$ws.Range("H24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H24").Value = 14.17004048583
$ws.Range("I24").NumberFormat = '#,##0'
$ws.Range("I24").Value = 3295
$ws.Range("J24").NumberFormat = '#,##0'
$ws.Range("J24").Value = 2224
$ws.Range("K24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K24").Value = 48.156474820143
$ws.Range("L24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L24").Value = 88.070776255707
$ws.Range("M24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M24").Value = -29.518716577540
$ws.Range("C25").NumberFormat = '#,##0'
$ws.Range("C25").Value = 23
$ws.Range("D25").NumberFormat = '#,##0'
$ws.Range("D25").Value = 12
$ws.Range("E25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E25").Value = 91.666666666666
$ws.Range("F25").NumberFormat = '#,##0'
$ws.Range("F25").Value = 88
$ws.Range("G25").NumberFormat = '#,##0'
$ws.Range("G25").Value = 49
$ws.Range("H25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H25").Value = 79.591836734693
$ws.Range("I25").NumberFormat = '#,##0'
$ws.Range("I25").Value = 893
$ws.Range("J25").NumberFormat = '#,##0'
$ws.Range("J25").Value = 838
$ws.Range("K25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K25").Value = 6.563245823389
$ws.Range("L25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L25").Value = 75.787401574803
$ws.Range("M25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M25").Value = 51.612903225806
$ws.Range("C23").Copy($ws.Range("D26"))
$ws.Range("E23").Copy($ws.Range("E26"))
$ws.Range("C23").Copy($ws.Range("F26"))
$ws.Range("H26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H26").Value = -100
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("C27").Value = 3
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 3
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = 0
$ws.Range("F27").NumberFormat = '#,##0'
$ws.Range("F27").Value = 20
$ws.Range("G27").NumberFormat = '#,##0'
$ws.Range("G27").Value = 11
$ws.Range("H27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H27").Value = 81.818181818181
$ws.Range("I27").NumberFormat = '#,##0'
$ws.Range("I27").Value = 219
$ws.Range("J27").NumberFormat = '#,##0'
$ws.Range("J27").Value = 149
$ws.Range("K27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K27").Value = 46.979865771812
$ws.Range("L27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L27").Value = 88.793103448275
$ws.Range("C23").Copy($ws.Range("G28"))
$ws.Range("E23").Copy($ws.Range("H28"))
$ws.Range("C23").Copy($ws.Range("G29"))
$ws.Range("E23").Copy($ws.Range("H29"))
$ws.Range("C30").NumberFormat = '#,##0'
$ws.Range("C30").Value = 1
$ws.Range("C23").Copy($ws.Range("D30"))
$ws.Range("E23").Copy($ws.Range("E30"))
$ws.Range("F30").NumberFormat = '#,##0'
$ws.Range("F30").Value = 1
$ws.Range("H30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H30").Value = 0
$ws.Range("I30").NumberFormat = '#,##0'
$ws.Range("I30").Value = 25
$ws.Range("K30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K30").Value = -30.555555555555
$ws.Range("L30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L30").Value = 733.333333333333
